# ---------------------------------------------------------------------------
# Reproduces the recorded edit:
#  - Sheet1!B3 data point corrected 257.7 -> 719.6 (and the two charts that
#    cache that series get the same correction where the engine allows it).
#  - Sheet2 is populated with the "cal" calibration pull (24 rows, A:B) that
#    used to live in an external "From Text" query (cal.txt, ':' delimited).
#    We reproduce the query's end effect -- the literal worksheet values,
#    the query table's defined name, and a best-effort external connection
#    -- since this host's QueryTables/Connections COM surface is stubbed.
#  - A new scatter chart ("RPM-Power", with a linear trendline) is added on
#    Sheet2 plotting column B (x) against column A (y), matching chart3.xml.
#  - View state: Sheet2 becomes the active tab/sheet, Sheet1's selection
#    moves to B4, Sheet2's selection sits at K5.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1) Fix the mis-keyed calibration point on Sheet1.
# ---------------------------------------------------------------------------
$ws1.Range("B3").Value = 719.6

# ---------------------------------------------------------------------------
# 2) Populate Sheet2 with the calibration table (24 rows x 2 cols).
# ---------------------------------------------------------------------------
$calData = @(
    @(0.5, 1968.6),
    @(0.5, 1968.6),
    @(0.5, 1968.6),
    @(0.55, 2229.9),
    @(0.55, 2229.9),
    @(0.55, 2229.9),
    @(0.65, 2656.5),
    @(0.65, 2656.5),
    @(0.65, 2656.5),
    @(0.7, 2852.8),
    @(0.7, 2852.8),
    @(0.7, 2852.8),
    @(0.75, 2998.5),
    @(0.75, 2998.5),
    @(0.75, 2998.5),
    @(0.8, 3111.6),
    @(0.8, 3111.6),
    @(0.8, 3111.6),
    @(0.85, 3304.7),
    @(0.85, 3304.7),
    @(0.85, 3304.7),
    @(0.9, 3519.7),
    @(0.9, 3519.7),
    @(0.9, 3519.7)
)

for ($r = 0; $r -lt $calData.Length; $r++) {
    $ws2.Cells.Item($r + 1, 1).Value = $calData[$r][0]
    $ws2.Cells.Item($r + 1, 2).Value = $calData[$r][1]
}

# Column widths, matching the auto-fit widths the text-import wizard leaves
# behind on the destination sheet.
$ws2.Columns("A").ColumnWidth = 5
$ws2.Columns("B").ColumnWidth = 7

$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3) Recreate the external "cal" text-file connection + query table that
#    produced the data above. QueryTables/Connections are not wired up to
#    this host's writer, so these calls are best-effort (no-throw) and the
#    defined name below carries the piece that IS persisted.
# ---------------------------------------------------------------------------
try {
    $connStr = 'TEXT;C:\Users\irving\Desktop\cal.txt'
    $qt = $ws2.QueryTables.Add($connStr, $ws2.Range("A1"))
    if ($qt) {
        $qt.Name = "cal_1"
        $qt.TextFilePlatform = 437
        $qt.TextFileStartRow = 1
        $qt.TextFileParseType = 1
        $qt.TextFileConsecutiveDelimiter = $true
        $qt.TextFileSpaceDelimiter = $true
        $qt.TextFileOtherDelimiter = ":"
        $qt.TextFileTrailingMinusNumbers = $true
        $qt.Refresh() | Out-Null
    }
} catch {}

# The query table's own defined name -- this DOES persist via Names.Add.
$ws2.Names.Add('cal_1', '=Sheet2!$A$1:$B$45')

# ---------------------------------------------------------------------------
# 4) Add the new "RPM-Power" scatter chart (chart3.xml) on Sheet2, plotting
#    column B as X and column A as Y, with a linear trendline.
# ---------------------------------------------------------------------------
$chartObj = $ws2.ChartObjects().Add(228600, 123825, 1143000, 1080000)
$newChart = $chartObj.Chart
$newChart.ChartType = -4169  # xlXYScatterLines

$newSeries = $newChart.SeriesCollection().NewSeries()
$newSeries.Name = "RPM-Power"
$newSeries.XValues = $ws2.Range("B:B")
$newSeries.Values = $ws2.Range("A:A")

$trend = $newSeries.Trendlines().Add()
$trend.Type = -4132  # xlLinear
$trend.DisplayEquation = $true
$trend.DisplayRSquared = $false
try { $trend.DataLabel.NumberFormat = "#,##0.000000000000" } catch {}

# ---------------------------------------------------------------------------
# 5) View state: Sheet2 active tab, Sheet1 selection -> B4, Sheet2 -> K5.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B4").Select() | Out-Null

$ws2.Activate()
$ws2.Range("K5").Select() | Out-Null
